# Updated coin price/volume/hour data per refreshed symbol list (run Sat Feb  4 05:11:08 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. Columns D (Price), E (Volume 1h) and G (Hora)
# hold numeric-looking text, so we force a text NumberFormat before assigning the
# Value, otherwise Excel would silently reinterpret them as numbers/percentages.
$updates = [ordered]@{
    "D2" = "329.14"
    "E2" = "1.92%"
    "G2" = "5"
    "D3" = "41.15"
    "E3" = "3.91%"
    "G3" = "5"
    "D4" = "5.617"
    "E4" = "-4.39%"
    "G4" = "5"
    "D5" = "0.08168"
    "E5" = "1.76%"
    "G5" = "5"
    "D6" = "2.048"
    "E6" = "6.02%"
    "G6" = "5"
    "D7" = "8.745"
    "E7" = "0.97%"
    "G7" = "5"
    "D8" = "4.540"
    "E8" = "-0.75%"
    "G8" = "5"
    "D9" = "2.931"
    "E9" = "-0.41%"
    "G9" = "5"
    "D10" = "0.9178"
    "E10" = "-1.48%"
    "G10" = "5"
    "D11" = "0.1254"
    "E11" = "-1.05%"
    "G11" = "5"
    "D12" = "0.1959"
    "E12" = "-0.02%"
    "G12" = "5"
    "D13" = "0.09349"
    "E13" = "1.58%"
    "G13" = "5"
    "D14" = "0.03689"
    "E14" = "4.34%"
    "G14" = "5"
    "E15" = "10.45%"
    "G15" = "5"
    "D16" = "0.001295"
    "E16" = "0.16%"
    "G16" = "5"
    "D17" = "0.006135"
    "E17" = "0.90%"
    "G17" = "5"
    "D18" = "3.434"
    "E18" = "2.62%"
    "G18" = "5"
    "E19" = "-2.12%"
    "G19" = "5"
    "D20" = "8.269"
    "E20" = "-5.29%"
    "G20" = "5"
    "D21" = "0.1393"
    "E21" = "-1.78%"
    "G21" = "5"
    "D22" = "0.2652"
    "E22" = "10.10%"
    "G22" = "5"
    "D23" = "0.04428"
    "E23" = "0.43%"
    "G23" = "5"
    "D24" = "0.001267"
    "E24" = "0.50%"
    "G24" = "5"
    "D25" = "0.004291"
    "E25" = "-2.43%"
    "G25" = "5"
    "E26" = "3.69%"
    "G26" = "5"
    "G27" = "5"
    "G28" = "5"
    "G29" = "5"
    "G30" = "5"
    "G31" = "5"
    "G32" = "5"
    "G33" = "5"
    "G34" = "5"
    "G35" = "5"
    "G36" = "5"
    "G37" = "5"
    "G38" = "5"
    "D39" = "0.02767"
    "E39" = "13.84%"
    "G39" = "5"
    "D40" = "0.05441"
    "E40" = "3.92%"
    "G40" = "5"
    "D41" = "0.007668"
    "E41" = "3.25%"
    "G41" = "5"
    "D42" = "0.009475"
    "E42" = "0.40%"
    "G42" = "5"
    "D43" = "0.1415"
    "E43" = "0.67%"
    "G43" = "5"
    "E44" = "-0.30%"
    "G44" = "5"
    "D45" = "0.01160"
    "E45" = "3.81%"
    "G45" = "5"
    "D46" = "0.00006898"
    "E46" = "2.39%"
    "G46" = "5"
    "E47" = "0.16%"
    "G47" = "5"
    "B48" = "BOLO"
    "C48" = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
    "D48" = "0.003544"
    "E48" = "18.07%"
    "G48" = "5"
    "B49" = "CoinbaseStockToken"
    "C49" = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
    "D49" = "0.002283"
    "E49" = "60.50%"
    "G49" = "5"
    "E50" = "0.16%"
    "G50" = "5"
    "D51" = "0.0002004"
    "E51" = "0.16%"
    "G51" = "5"
}

$textColumns = @("D", "E", "G")

foreach ($ref in $updates.Keys) {
    $col = $ref -replace "[0-9]+$", ""
    $range = $ws.Range($ref)
    if ($textColumns -contains $col) {
        $range.NumberFormat = "@"
    }
    $range.Value = $updates[$ref]
}
